$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-367): update date value from 45184 to 45186
$ws.Range("C2:C367").Value = 45186

# Columns S, T, V, W, X, Y (rows 2-44): add the record name as the
# second (friendly-name) argument of the HYPERLINK formula.
$cols = @("S", "T", "V", "W", "X", "Y")
for ($row = 2; $row -le 44; $row++) {
    $name = $ws.Range("A$row").Value()
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$row")
        $formula = $cell.Formula()
        if ($formula -ne $null -and $formula -match 'HYPERLINK\("([^"]*)"\)') {
            $url = $matches[1]
            $newFormula = '=HYPERLINK("' + $url + '", "' + $name + '")'
            $cell.Formula = $newFormula
        }
    }
}
